$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "kesavan@congruentglobal.com"
$ws.Range("A3").Value = "sathish.j@congruentglobal.com"

$ws.Range("E2").Value = "Kesavan"
$ws.Range("F2").Value = "R"

$ws.Range("E3").Value = "Sathish"
$ws.Range("F3").Value = "J"

$ws.Range("D4").Select()
